$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.130.35"
$ws.Range("E2").Value = "  +0.98%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.540.71"
$ws.Range("E3").Value = "  +1.07%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "613.79"
$ws.Range("E5").Value = "  +1.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.39"
$ws.Range("E6").Value = "  +2.83%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.541.87"
$ws.Range("E7").Value = "  +1.11%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.490"
$ws.Range("E9").Value = "  +2.01%  "
$ws.Range("E10").Value = "  +4.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.65"
$ws.Range("E11").Value = "  +9.59%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.436"
$ws.Range("E12").Value = "  +3.55%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "33.01"
$ws.Range("E13").Value = "  +5.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000219"
$ws.Range("E14").Value = "  +1.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.138.25"
$ws.Range("E15").Value = "  +1.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.538.08"
$ws.Range("E16").Value = "  +0.74%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.043.81"
$ws.Range("E17").Value = "  +0.69%  "
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.65"
$ws.Range("E19").Value = "  +3.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.79"
$ws.Range("E20").Value = "  +4.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.90"
$ws.Range("E21").Value = "  +9.50%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "452.56"
$ws.Range("E22").Value = "  +1.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.640"
$ws.Range("E23").Value = "  +3.39%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.22"
$ws.Range("E24").Value = "  +1.14%  "
$ws.Range("E25").Value = "  +3.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.678.85"
$ws.Range("E26").Value = "  +0.99%  "
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.11"
$ws.Range("E28").Value = "  +11.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.25"
$ws.Range("E29").Value = "  +1.63%  "
$ws.Range("E30").Value = "  +10.46%  "
$ws.Range("E31").Value = "  +2.47%  "
$ws.Range("E32").Value = "  +4.38%  "
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.98"
$ws.Range("E34").Value = "  +1.33%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.28"
$ws.Range("E35").Value = "  +3.78%  "
$ws.Range("E36").Value = "  +3.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.532.17"
$ws.Range("E37").Value = "  +1.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.14"
$ws.Range("E38").Value = "  +1.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.35"
$ws.Range("E40").Value = "  +8.56%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  -0.10%  "
$ws.Range("B42").Value = "Hedera"
$ws.Range("C42").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0910"
$ws.Range("E42").Value = "  +3.64%  "
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "175.07"
$ws.Range("E43").Value = "  -1.91%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.58"
$ws.Range("E44").Value = "  +3.89%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "30.96"
$ws.Range("E45").Value = "  +10.81%  "
$ws.Range("E46").Value = "  +1.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "46.41"
$ws.Range("E47").Value = "  +2.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.34"
$ws.Range("E48").Value = "  +8.72%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.58"
$ws.Range("E49").Value = "  +1.51%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.74"
$ws.Range("E50").Value = "  +2.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.258"
$ws.Range("E51").Value = "  +6.32%  "
